# "Electricity Technology Shareweights.xlsx" update (commit: "update becf, etle, ets")
#
# On the "ETS" sheet:
#   - Row 6 ("onshore wind") shareweights for every year column (B:AF,
#     i.e. 2020-2050) go from 1 to 3.
#   - Row 7 ("solar pv") shareweights for every year column (B:AF) go
#     from 1 to 2.
#
# The workbook is also left with the "ETS" sheet active/selected at cell
# A6 (it previously had the "About" sheet active and ETS's own selection
# sitting on B17).

$wb = $excel.ActiveWorkbook
$etsSheet = $wb.Worksheets.Item("ETS")

# Onshore wind (row 6): 1 -> 3 for 2020 (col B) .. 2050 (col AF)
$etsSheet.Range("B6:AF6").Value = 3

# Solar PV (row 7): 1 -> 2 for 2020 (col B) .. 2050 (col AF)
$etsSheet.Range("B7:AF7").Value = 2

# Make "ETS" the active sheet with A6 as the selected cell.
$etsSheet.Activate()
$etsSheet.Range("A6").Select()
